$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Complete the previously truncated text in G33 ("No error and branch gets " -> "...deleted")
$ws.Range("G33").Value = "No error and branch gets deleted"

# Clone formatting from the last existing data row (row 33) onto the two new rows
# so the new cells pick up the same styles (center-aligned numbers in column E,
# left-aligned text in columns F/G) without introducing new style entries.
$ws.Range("E33").Copy()
$ws.Range("E34:E35").PasteSpecial(-4122)

$ws.Range("F33").Copy()
$ws.Range("F34:F35").PasteSpecial(-4122)
$ws.Range("G34:G35").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# Populate the new rows, column by column:
# Row 34: git branch -d origin develop  -> Deletes develop branch from github website
# Row 35: git push origin develop       -> pushes develop branch to github
$ws.Range("E34").Value = 32
$ws.Range("E35").Value = 33

$ws.Range("F34").Value = "git branch -d origin develop"
$ws.Range("F35").Value = "git push origin develop"

$ws.Range("G34").Value = "Deletes develop branch from github website"
$ws.Range("G35").Value = "pushes develop branch to github"

# Mirror the author's final selection/scroll position
$ws.Range("G35").Select()
